{"js": "// Update the date line and the 25 division-problem answers in the table\n// to the values from the next day's worksheet (2024-11-13 Wednesday).\nconst replacements = [\n  [\"2024-11-12 Tuesday\", \"2024-11-13 Wednesday\"],\n  [\"10\u00f77=1, 3\", \"73\u00f78=9, 1\"],\n  [\"13\u00f77=1, 6\", \"23\u00f78=2, 7\"],\n  [\"99\u00f75=19, 4\", \"10\u00f75=2, 0\"],\n  [\"11\u00f79=1, 2\", \"31\u00f73=10, 1\"],\n  [\"79\u00f72=39, 1\", \"67\u00f77=9, 4\"],\n  [\"33\u00f79=3, 6\", \"30\u00f73=10, 0\"],\n  [\"18\u00f78=2, 2\", \"28\u00f72=14, 0\"],\n  [\"26\u00f75=5, 1\", \"65\u00f77=9, 2\"],\n  [\"63\u00f74=15, 3\", \"89\u00f77=12, 5\"],\n  [\"71\u00f79=7, 8\", \"29\u00f75=5, 4\"],\n  [\"51\u00f75=10, 1\", \"29\u00f73=9, 2\"],\n  [\"52\u00f75=10, 2\", \"90\u00f78=11, 2\"],\n  [\"16\u00f77=2, 2\", \"43\u00f78=5, 3\"],\n  [\"71\u00f73=23, 2\", \"77\u00f77=11, 0\"],\n  [\"56\u00f75=11, 1\", \"98\u00f78=12, 2\"],\n  [\"13\u00f79=1, 4\", \"97\u00f77=13, 6\"],\n  [\"94\u00f76=15, 4\", \"74\u00f75=14, 4\"],\n  [\"27\u00f72=13, 1\", \"18\u00f73=6, 0\"],\n  [\"61\u00f77=8, 5\", \"68\u00f75=13, 3\"],\n  [\"30\u00f74=7, 2\", \"10\u00f76=1, 4\"],\n  [\"32\u00f73=10, 2\", \"44\u00f75=8, 4\"],\n  [\"56\u00f79=6, 2\", \"41\u00f72=20, 1\"],\n  [\"89\u00f78=11, 1\", \"74\u00f78=9, 2\"],\n  [\"26\u00f79=2, 8\", \"81\u00f78=10, 1\"],\n  [\"60\u00f76=10, 0\", \"48\u00f75=9, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division-problem answers in the table\n# to the values from the next day's worksheet (2024-11-13 Wednesday).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-11-12 Tuesday\", \"2024-11-13 Wednesday\"),\n    @(\"10\u00f77=1, 3\", \"73\u00f78=9, 1\"),\n    @(\"13\u00f77=1, 6\", \"23\u00f78=2, 7\"),\n    @(\"99\u00f75=19, 4\", \"10\u00f75=2, 0\"),\n    @(\"11\u00f79=1, 2\", \"31\u00f73=10, 1\"),\n    @(\"79\u00f72=39, 1\", \"67\u00f77=9, 4\"),\n    @(\"33\u00f79=3, 6\", \"30\u00f73=10, 0\"),\n    @(\"18\u00f78=2, 2\", \"28\u00f72=14, 0\"),\n    @(\"26\u00f75=5, 1\", \"65\u00f77=9, 2\"),\n    @(\"63\u00f74=15, 3\", \"89\u00f77=12, 5\"),\n    @(\"71\u00f79=7, 8\", \"29\u00f75=5, 4\"),\n    @(\"51\u00f75=10, 1\", \"29\u00f73=9, 2\"),\n    @(\"52\u00f75=10, 2\", \"90\u00f78=11, 2\"),\n    @(\"16\u00f77=2, 2\", \"43\u00f78=5, 3\"),\n    @(\"71\u00f73=23, 2\", \"77\u00f77=11, 0\"),\n    @(\"56\u00f75=11, 1\", \"98\u00f78=12, 2\"),\n    @(\"13\u00f79=1, 4\", \"97\u00f77=13, 6\"),\n    @(\"94\u00f76=15, 4\", \"74\u00f75=14, 4\"),\n    @(\"27\u00f72=13, 1\", \"18\u00f73=6, 0\"),\n    @(\"61\u00f77=8, 5\", \"68\u00f75=13, 3\"),\n    @(\"30\u00f74=7, 2\", \"10\u00f76=1, 4\"),\n    @(\"32\u00f73=10, 2\", \"44\u00f75=8, 4\"),\n    @(\"56\u00f79=6, 2\", \"41\u00f72=20, 1\"),\n    @(\"89\u00f78=11, 1\", \"74\u00f78=9, 2\"),\n    @(\"26\u00f79=2, 8\", \"81\u00f78=10, 1\"),\n    @(\"60\u00f76=10, 0\", \"48\u00f75=9, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $f = $d.Content.Find\n    $f.ClearFormatting()\n    $f.Replacement.ClearFormatting()\n    $f.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
